$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update totals at the top of the statement ---
# "VALOR MORA" total
$ws.Range("E11").Value = 264000
# "Cant. Periodos" (number of periods)
$ws.Range("F13").Value = 2

# --- Insert two new rows into the detail table (it grows from 2 to 4 rows) ---
# Existing data table occupies rows 16 (normal) and 17 (bottom-border / last row).
# We insert 2 rows right after row 16, which pushes the old row 17 down to row 19,
# and copy row 16's formatting into the two freshly inserted rows so every row in
# the table keeps the correct borders/number formats.
$ws.Rows("17:18").Insert()
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$ws.Range("B18:J18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Populate the 4 detail rows ---
# Row 16: JANER ENRIQUE CHAMORRO ROCHA, period 2507
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1063495346"
$ws.Range("D16").Value = "JANER ENRIQUE CHAMORRO ROCHA"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 80000
$ws.Range("G16").Value = 2000000

# Row 17: JANER ENRIQUE CHAMORRO ROCHA, period 2506
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1063495346"
$ws.Range("D17").Value = "JANER ENRIQUE CHAMORRO ROCHA"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 66667
$ws.Range("G17").Value = 2000000

# Row 18: NEIVER DE JESUS CHAMORRO ROCHA, period 2507
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1124503428"
$ws.Range("D18").Value = "NEIVER DE JESUS CHAMORRO ROCHA"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 64000
$ws.Range("G18").Value = 1600000

# Row 19: NEIVER DE JESUS CHAMORRO ROCHA, period 2506
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1124503428"
$ws.Range("D19").Value = "NEIVER DE JESUS CHAMORRO ROCHA"
$ws.Range("E19").Value = "2506"
$ws.Range("F19").Value = 53333
$ws.Range("G19").Value = 1600000

# --- Re-fit the table columns now that some cells hold wider content ---
# (ColumnWidth is in "characters"; OOXML <col width> = ColumnWidth + 5/6, quantized
# to 1/6ths here, so these are the closest achievable values to the real best-fit
# widths Excel would compute for the new content.)
$ws.Columns("B").ColumnWidth = 17.666666666666668
$ws.Columns("C").ColumnWidth = 15.833333333333334
$ws.Columns("E").ColumnWidth = 12.666666666666666
$ws.Columns("F").ColumnWidth = 9.333333333333332
$ws.Columns("G").ColumnWidth = 13.5
$ws.Columns("H").ColumnWidth = 18.5
$ws.Columns("I").ColumnWidth = 17.333333333333336
$ws.Columns("J").ColumnWidth = 14.166666666666666
